$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 1733000000.0
$ws.Range("D4").Value = 1780000000.0
$ws.Range("E4").Value = 1695000000.0
$ws.Range("F4").Value = 1615000000.0
$ws.Range("G4").Value = 1655000000.0

$ws.Range("C15").Value = 779000000.0
$ws.Range("D15").Value = 836000000.0
$ws.Range("E15").Value = 770000000.0
$ws.Range("F15").Value = 863000000.0
$ws.Range("G15").Value = 719000000.0

$ws.Range("C22").Value = 57000000.0
$ws.Range("D22").Value = 101000000.0
$ws.Range("E22").Value = 93000000.0
$ws.Range("F22").Value = 89000000.0
$ws.Range("G22").Value = 67000000.0

$ws.Range("B26").Value = 1248000000.0
$ws.Range("B27").Value = 55000000.0
$ws.Range("B28").Value = 8948000000.0
$ws.Range("B29").Value = 8354000000.0
$ws.Range("B30").Value = 2080000000.0

$ws.Range("C33").Value = 52524000.0

$ws.Range("B35").Value = 1818000000.0
$ws.Range("B36").Value = 2380000000.0
